$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 178.83
$ws.Range("I15").Value = 178.83
$ws.Range("K15").Value = 536.49
$ws.Range("M15").Value = -367.49
$ws.Range("H19").Value = 1100.5
$ws.Range("I19").Value = 2259.2
$ws.Range("J19").Value = 376.3125
$ws.Range("K19").Value = 2259.2
$ws.Range("L19").Value = 376.3125
$ws.Range("M19").Value = -2084.2
$ws.Range("N19").Value = -726.3125
$ws.Range("H28").Value = 886.5333000000001
$ws.Range("I28").Value = 874.0625
$ws.Range("J28").Value = 900.7857
$ws.Range("K28").Value = 874.0625
$ws.Range("L28").Value = 900.7857
$ws.Range("M28").Value = -389.0625
$ws.Range("N28").Value = -1870.7857
$ws.Range("H112").Value = 55556884
$ws.Range("I112").Value = 699.6
$ws.Range("J112").Value = 76924650
$ws.Range("K112").Value = 2098.8
$ws.Range("L112").Value = 230773950
$ws.Range("M112").Value = -990.8000000000002
$ws.Range("N112").Value = -230776166
$ws.Range("H116").Value = 2100.3572
$ws.Range("I116").Value = 2088.125
$ws.Range("J116").Value = 2116.6667
$ws.Range("K116").Value = 2088.125
$ws.Range("L116").Value = 2116.6667
$ws.Range("M116").Value = 1353.875
$ws.Range("N116").Value = -9000.6667
$ws.Range("H125").Value = 1185.9584
$ws.Range("I125").Value = 788.8
$ws.Range("J125").Value = 1847.8889
$ws.Range("K125").Value = 7099.2
$ws.Range("L125").Value = 16631.0001
$ws.Range("M125").Value = -4639.2
$ws.Range("N125").Value = -21551.0001
$ws.Range("H129").Value = 960.04
$ws.Range("J129").Value = 1160.9736
$ws.Range("L129").Value = 3482.9208
$ws.Range("N129").Value = -13482.9208
$ws.Range("H138").Value = 4051.8118
$ws.Range("I138").Value = 3258.6333
$ws.Range("J138").Value = 4484.4546
$ws.Range("K138").Value = 9775.8999
$ws.Range("L138").Value = 13453.3638
$ws.Range("M138").Value = -4635.8999
$ws.Range("N138").Value = -23733.3638
$ws.Range("H141").Value = 2217.5
$ws.Range("I141").Value = 1440.2941
$ws.Range("J141").Value = 4860
$ws.Range("K141").Value = 4320.8823
$ws.Range("L141").Value = 14580
$ws.Range("M141").Value = 859.1176999999998
$ws.Range("N141").Value = -24940

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23185.154
$ws.Range("I32").Value = 16747.51
$ws.Range("J32").Value = 128333.336
$ws.Range("K32").Value = 16747.51
$ws.Range("L32").Value = 128333.336
$ws.Range("M32").Value = -16460.51
$ws.Range("N32").Value = -128907.336
$ws.Range("H122").Value = 1372.1904
$ws.Range("I122").Value = 1009.7143
$ws.Range("J122").Value = 1734.6666
$ws.Range("K122").Value = 3029.1429
$ws.Range("L122").Value = 5203.9998
$ws.Range("M122").Value = -579.1428999999998
$ws.Range("N122").Value = -10103.9998
$ws.Range("H132").Value = 909426.7
$ws.Range("I132").Value = 1114327.6
$ws.Range("J132").Value = 2008.2858
$ws.Range("K132").Value = 3342982.8
$ws.Range("L132").Value = 6024.857400000001
$ws.Range("M132").Value = -3340452.8
$ws.Range("N132").Value = -11084.8574

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 462.7143
$ws.Range("I94").Value = 373.16666
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 373.16666
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 77.83334000000002
$ws.Range("N94").Value = -1902
$ws.Range("H107").Value = 2019355.1
$ws.Range("I107").Value = 7047743.5
$ws.Range("J107").Value = 7999.8
$ws.Range("K107").Value = 7047743.5
$ws.Range("L107").Value = 7999.8
$ws.Range("M107").Value = -7045823.5
$ws.Range("N107").Value = -11839.8
$ws.Range("H134").Value = 5129739.5
$ws.Range("I134").Value = 5715824.5
$ws.Range("J134").Value = 1497.5
$ws.Range("K134").Value = 17147473.5
$ws.Range("L134").Value = 4492.5
$ws.Range("M134").Value = -17144938.5
$ws.Range("N134").Value = -9562.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 18463.857
$ws.Range("J68").Value = 19163.166
$ws.Range("L68").Value = 19163.166
$ws.Range("N68").Value = -20661.166
$ws.Range("H71").Value = 18463.857
$ws.Range("J71").Value = 19163.166
$ws.Range("L71").Value = 57489.49800000001
$ws.Range("N71").Value = -64977.49800000001
$ws.Range("H122").Value = 8184.3335
$ws.Range("I122").Value = 14131.5
$ws.Range("J122").Value = 1387.5714
$ws.Range("K122").Value = 42394.5
$ws.Range("L122").Value = 4162.7142
$ws.Range("M122").Value = -39944.5
$ws.Range("N122").Value = -9062.7142

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 35357544
$ws.Range("J131").Value = 17243872
$ws.Range("L131").Value = 51731616
$ws.Range("N131").Value = -51741696

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = $null
$ws.Range("N28").Value = 0
$ws.Range("H80").Value = 57361.05
$ws.Range("I80").Value = 2425.9092
$ws.Range("J80").Value = 124504
$ws.Range("K80").Value = 2425.9092
$ws.Range("L80").Value = 124504
$ws.Range("M80").Value = -1427.9092
$ws.Range("N80").Value = -126500
$ws.Range("H83").Value = 57361.05
$ws.Range("I83").Value = 2425.9092
$ws.Range("J83").Value = 124504
$ws.Range("K83").Value = 12129.546
$ws.Range("L83").Value = 622520
$ws.Range("M83").Value = -7137.546
$ws.Range("N83").Value = -632504
$ws.Range("H126").Value = 1191.3334
$ws.Range("I126").Value = 796.25
$ws.Range("K126").Value = 2388.75
$ws.Range("M126").Value = 81.25
$ws.Range("H132").Value = 1649.1072
$ws.Range("I132").Value = 1198.15
$ws.Range("J132").Value = 2776.5
$ws.Range("K132").Value = 3594.45
$ws.Range("L132").Value = 8329.5
$ws.Range("M132").Value = -1064.45
$ws.Range("N132").Value = -13389.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 359.73334
$ws.Range("I22").Value = 318.18182
$ws.Range("J22").Value = 474
$ws.Range("K22").Value = 318.18182
$ws.Range("L22").Value = 474
$ws.Range("M22").Value = -23.18182000000002
$ws.Range("N22").Value = -1064
$ws.Range("H27").Value = 359.73334
$ws.Range("I27").Value = 318.18182
$ws.Range("J27").Value = 474
$ws.Range("K27").Value = 318.18182
$ws.Range("L27").Value = 474
$ws.Range("M27").Value = -211.18182
$ws.Range("N27").Value = -688
$ws.Range("H55").Value = 112.789474
$ws.Range("I55").Value = 84.5
$ws.Range("J55").Value = 144.22223
$ws.Range("K55").Value = 84.5
$ws.Range("L55").Value = 144.22223
$ws.Range("M55").Value = 88.5
$ws.Range("N55").Value = -490.22223
$ws.Range("H82").Value = 1586.9131
$ws.Range("I82").Value = 1554.5454
$ws.Range("K82").Value = 1554.5454
$ws.Range("M82").Value = -1193.5454
$ws.Range("H85").Value = 1586.9131
$ws.Range("I85").Value = 1554.5454
$ws.Range("K85").Value = 1554.5454
$ws.Range("M85").Value = -306.5454
$ws.Range("H122").Value = 1950.5769
$ws.Range("I122").Value = 1841.4615
$ws.Range("J122").Value = 2277.923
$ws.Range("K122").Value = 5524.3845
$ws.Range("L122").Value = 6833.768999999999
$ws.Range("M122").Value = -3074.3845
$ws.Range("N122").Value = -11733.769
